$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the now-unused trailing column (BA) and trailing rows (23:24)
$ws.Columns("BA:BA").Clear() | Out-Null
$ws.Rows("23:24").Clear() | Out-Null

# Build the full corrected data grid (rows 1-22, columns A-AZ)
$arr = New-Object 'object[,]' 22,52
$arr[0,1] = 39583
$arr[0,2] = 39765
$arr[0,3] = 39948
$arr[0,4] = 40130
$arr[0,5] = 40310
$arr[0,6] = 40494
$arr[0,7] = 40676
$arr[0,8] = 40862
$arr[0,9] = 41044
$arr[0,10] = 41228
$arr[0,11] = 41409
$arr[0,12] = 41592
$arr[0,13] = 41774
$arr[0,14] = 41957
$arr[0,15] = 42137
$arr[0,16] = 42321
$arr[0,17] = 42503
$arr[0,18] = 42689
$arr[0,19] = 42867
$arr[0,20] = 43053
$arr[0,21] = 43145
$arr[0,22] = 43235
$arr[0,23] = 43326
$arr[0,24] = 43418
$arr[0,25] = 43510
$arr[0,26] = 43600
$arr[0,27] = 43691
$arr[0,28] = 43783
$arr[0,29] = 43875
$arr[0,30] = 43966
$arr[0,31] = 44068
$arr[0,32] = 44159
$arr[0,33] = 44251
$arr[0,34] = 44341
$arr[0,35] = 44432
$arr[0,36] = 44525
$arr[0,37] = 44617
$arr[0,38] = 44706
$arr[0,39] = 44798
$arr[0,40] = 44890
$arr[0,41] = 44981
$arr[0,42] = 45071
$arr[0,43] = 45163
$arr[0,44] = 45254
$arr[0,45] = 45345
$arr[0,46] = 45436
$arr[0,47] = 45534
$arr[0,48] = 45618
$arr[0,49] = 45713
$arr[0,50] = 45800
$arr[0,51] = 45891
$arr[1,0] = 39813
$arr[2,0] = 40178
$arr[2,1] = 1.328558632615739
$arr[2,2] = 0.4282194198276246
$arr[2,3] = -1.435981453719049
$arr[2,4] = -1.61188520630966
$arr[2,5] = -1.61188520630966
$arr[2,6] = -1.61188520630966
$arr[2,7] = -1.61188520630966
$arr[2,8] = -1.61188520630966
$arr[2,9] = -1.61188520630966
$arr[2,10] = -1.61188520630966
$arr[2,11] = -1.61188520630966
$arr[2,12] = -1.61188520630966
$arr[2,13] = -1.61188520630966
$arr[2,14] = -1.61188520630966
$arr[2,15] = -1.61188520630966
$arr[2,16] = -1.61188520630966
$arr[2,17] = -1.61188520630966
$arr[2,18] = -1.61188520630966
$arr[2,19] = -1.61188520630966
$arr[2,20] = -1.61188520630966
$arr[2,21] = -1.61188520630966
$arr[2,22] = -1.61188520630966
$arr[2,23] = -1.61188520630966
$arr[2,24] = -1.61188520630966
$arr[2,25] = -1.61188520630966
$arr[2,26] = -1.61188520630966
$arr[2,27] = -1.61188520630966
$arr[2,28] = -1.61188520630966
$arr[2,29] = -1.61188520630966
$arr[2,30] = -1.61188520630966
$arr[2,31] = -1.61188520630966
$arr[2,32] = -1.61188520630966
$arr[2,33] = -1.61188520630966
$arr[2,34] = -1.61188520630966
$arr[2,35] = -1.61188520630966
$arr[2,36] = -1.61188520630966
$arr[2,37] = -1.61188520630966
$arr[2,38] = -1.61188520630966
$arr[2,39] = -1.61188520630966
$arr[2,40] = -1.61188520630966
$arr[2,41] = -1.61188520630966
$arr[2,42] = -1.61188520630966
$arr[2,43] = -1.61188520630966
$arr[2,44] = -1.61188520630966
$arr[2,45] = -1.61188520630966
$arr[2,46] = -1.61188520630966
$arr[2,47] = -1.61188520630966
$arr[2,48] = -1.61188520630966
$arr[2,49] = -1.61188520630966
$arr[2,50] = -1.61188520630966
$arr[2,51] = -1.61188520630966
$arr[3,0] = 40543
$arr[3,1] = 1.356849765318358
$arr[3,2] = 0.5762410323606026
$arr[3,3] = -0.7704417043119083
$arr[3,4] = -0.8235211753995442
$arr[3,5] = 0.406633294022174
$arr[3,6] = 1.128600547465064
$arr[3,7] = 1.128600547465064
$arr[3,8] = 1.128600547465064
$arr[3,9] = 1.128600547465064
$arr[3,10] = 1.128600547465064
$arr[3,11] = 1.128600547465064
$arr[3,12] = 1.128600547465064
$arr[3,13] = 1.128600547465064
$arr[3,14] = 1.128600547465064
$arr[3,15] = 1.128600547465064
$arr[3,16] = 1.128600547465064
$arr[3,17] = 1.128600547465064
$arr[3,18] = 1.128600547465064
$arr[3,19] = 1.128600547465064
$arr[3,20] = 1.128600547465064
$arr[3,21] = 1.128600547465064
$arr[3,22] = 1.128600547465064
$arr[3,23] = 1.128600547465064
$arr[3,24] = 1.128600547465064
$arr[3,25] = 1.128600547465064
$arr[3,26] = 1.128600547465064
$arr[3,27] = 1.128600547465064
$arr[3,28] = 1.128600547465064
$arr[3,29] = 1.128600547465064
$arr[3,30] = 1.128600547465064
$arr[3,31] = 1.128600547465064
$arr[3,32] = 1.128600547465064
$arr[3,33] = 1.128600547465064
$arr[3,34] = 1.128600547465064
$arr[3,35] = 1.128600547465064
$arr[3,36] = 1.128600547465064
$arr[3,37] = 1.128600547465064
$arr[3,38] = 1.128600547465064
$arr[3,39] = 1.128600547465064
$arr[3,40] = 1.128600547465064
$arr[3,41] = 1.128600547465064
$arr[3,42] = 1.128600547465064
$arr[3,43] = 1.128600547465064
$arr[3,44] = 1.128600547465064
$arr[3,45] = 1.128600547465064
$arr[3,46] = 1.128600547465064
$arr[3,47] = 1.128600547465064
$arr[3,48] = 1.128600547465064
$arr[3,49] = 1.128600547465064
$arr[3,50] = 1.128600547465064
$arr[3,51] = 1.128600547465064
$arr[4,0] = 40908
$arr[4,3] = -0.6484173407089511
$arr[4,4] = -0.735462857474789
$arr[4,5] = -0.1034614224434405
$arr[4,6] = 1.121293995080253
$arr[4,7] = 2.185496833134781
$arr[4,8] = 2.397632938760519
$arr[4,9] = 2.397632938760519
$arr[4,10] = 2.397632938760519
$arr[4,11] = 2.397632938760519
$arr[4,12] = 2.397632938760519
$arr[4,13] = 2.397632938760519
$arr[4,14] = 2.397632938760519
$arr[4,15] = 2.397632938760519
$arr[4,16] = 2.397632938760519
$arr[4,17] = 2.397632938760519
$arr[4,18] = 2.397632938760519
$arr[4,19] = 2.397632938760519
$arr[4,20] = 2.397632938760519
$arr[4,21] = 2.397632938760519
$arr[4,22] = 2.397632938760519
$arr[4,23] = 2.397632938760519
$arr[4,24] = 2.397632938760519
$arr[4,25] = 2.397632938760519
$arr[4,26] = 2.397632938760519
$arr[4,27] = 2.397632938760519
$arr[4,28] = 2.397632938760519
$arr[4,29] = 2.397632938760519
$arr[4,30] = 2.397632938760519
$arr[4,31] = 2.397632938760519
$arr[4,32] = 2.397632938760519
$arr[4,33] = 2.397632938760519
$arr[4,34] = 2.397632938760519
$arr[4,35] = 2.397632938760519
$arr[4,36] = 2.397632938760519
$arr[4,37] = 2.397632938760519
$arr[4,38] = 2.397632938760519
$arr[4,39] = 2.397632938760519
$arr[4,40] = 2.397632938760519
$arr[4,41] = 2.397632938760519
$arr[4,42] = 2.397632938760519
$arr[4,43] = 2.397632938760519
$arr[4,44] = 2.397632938760519
$arr[4,45] = 2.397632938760519
$arr[4,46] = 2.397632938760519
$arr[4,47] = 2.397632938760519
$arr[4,48] = 2.397632938760519
$arr[4,49] = 2.397632938760519
$arr[4,50] = 2.397632938760519
$arr[4,51] = 2.397632938760519
$arr[5,0] = 41274
$arr[5,5] = -0.2218152533720597
$arr[5,6] = 0.1480821602630744
$arr[5,7] = 0.6652762968575532
$arr[5,8] = 1.665250327443002
$arr[5,9] = 0.8574941660507873
$arr[5,10] = 0.5991205513815823
$arr[5,11] = 0.5991205513815823
$arr[5,12] = 0.5991205513815823
$arr[5,13] = 0.5991205513815823
$arr[5,14] = 0.5991205513815823
$arr[5,15] = 0.5991205513815823
$arr[5,16] = 0.5991205513815823
$arr[5,17] = 0.5991205513815823
$arr[5,18] = 0.5991205513815823
$arr[5,19] = 0.5991205513815823
$arr[5,20] = 0.5991205513815823
$arr[5,21] = 0.5991205513815823
$arr[5,22] = 0.5991205513815823
$arr[5,23] = 0.5991205513815823
$arr[5,24] = 0.5991205513815823
$arr[5,25] = 0.5991205513815823
$arr[5,26] = 0.5991205513815823
$arr[5,27] = 0.5991205513815823
$arr[5,28] = 0.5991205513815823
$arr[5,29] = 0.5991205513815823
$arr[5,30] = 0.5991205513815823
$arr[5,31] = 0.5991205513815823
$arr[5,32] = 0.5991205513815823
$arr[5,33] = 0.5991205513815823
$arr[5,34] = 0.5991205513815823
$arr[5,35] = 0.5991205513815823
$arr[5,36] = 0.5991205513815823
$arr[5,37] = 0.5991205513815823
$arr[5,38] = 0.5991205513815823
$arr[5,39] = 0.5991205513815823
$arr[5,40] = 0.5991205513815823
$arr[5,41] = 0.5991205513815823
$arr[5,42] = 0.5991205513815823
$arr[5,43] = 0.5991205513815823
$arr[5,44] = 0.5991205513815823
$arr[5,45] = 0.5991205513815823
$arr[5,46] = 0.5991205513815823
$arr[5,47] = 0.5991205513815823
$arr[5,48] = 0.5991205513815823
$arr[5,49] = 0.5991205513815823
$arr[5,50] = 0.5991205513815823
$arr[5,51] = 0.5991205513815823
$arr[6,0] = 41639
$arr[6,7] = 0.6294804885647043
$arr[6,8] = 1.617722252039
$arr[6,9] = 1.693469135756587
$arr[6,10] = 1.079796209653616
$arr[6,11] = 0.1494732105682406
$arr[6,12] = 0.5555179840670776
$arr[6,13] = 0.5555179840670776
$arr[6,14] = 0.5555179840670776
$arr[6,15] = 0.5555179840670776
$arr[6,16] = 0.5555179840670776
$arr[6,17] = 0.5555179840670776
$arr[6,18] = 0.5555179840670776
$arr[6,19] = 0.5555179840670776
$arr[6,20] = 0.5555179840670776
$arr[6,21] = 0.5555179840670776
$arr[6,22] = 0.5555179840670776
$arr[6,23] = 0.5555179840670776
$arr[6,24] = 0.5555179840670776
$arr[6,25] = 0.5555179840670776
$arr[6,26] = 0.5555179840670776
$arr[6,27] = 0.5555179840670776
$arr[6,28] = 0.5555179840670776
$arr[6,29] = 0.5555179840670776
$arr[6,30] = 0.5555179840670776
$arr[6,31] = 0.5555179840670776
$arr[6,32] = 0.5555179840670776
$arr[6,33] = 0.5555179840670776
$arr[6,34] = 0.5555179840670776
$arr[6,35] = 0.5555179840670776
$arr[6,36] = 0.5555179840670776
$arr[6,37] = 0.5555179840670776
$arr[6,38] = 0.5555179840670776
$arr[6,39] = 0.5555179840670776
$arr[6,40] = 0.5555179840670776
$arr[6,41] = 0.5555179840670776
$arr[6,42] = 0.5555179840670776
$arr[6,43] = 0.5555179840670776
$arr[6,44] = 0.5555179840670776
$arr[6,45] = 0.5555179840670776
$arr[6,46] = 0.5555179840670776
$arr[6,47] = 0.5555179840670776
$arr[6,48] = 0.5555179840670776
$arr[6,49] = 0.5555179840670776
$arr[6,50] = 0.5555179840670776
$arr[6,51] = 0.5555179840670776
$arr[7,0] = 42004
$arr[7,9] = 1.77975641852226
$arr[7,10] = 1.488234279941625
$arr[7,11] = 0.8024032015999882
$arr[7,12] = 1.374377011838535
$arr[7,13] = 1.656936590801972
$arr[7,14] = 1.749602965204744
$arr[7,15] = 1.749602965204744
$arr[7,16] = 1.749602965204744
$arr[7,17] = 1.749602965204744
$arr[7,18] = 1.749602965204744
$arr[7,19] = 1.749602965204744
$arr[7,20] = 1.749602965204744
$arr[7,21] = 1.749602965204744
$arr[7,22] = 1.749602965204744
$arr[7,23] = 1.749602965204744
$arr[7,24] = 1.749602965204744
$arr[7,25] = 1.749602965204744
$arr[7,26] = 1.749602965204744
$arr[7,27] = 1.749602965204744
$arr[7,28] = 1.749602965204744
$arr[7,29] = 1.749602965204744
$arr[7,30] = 1.749602965204744
$arr[7,31] = 1.749602965204744
$arr[7,32] = 1.749602965204744
$arr[7,33] = 1.749602965204744
$arr[7,34] = 1.749602965204744
$arr[7,35] = 1.749602965204744
$arr[7,36] = 1.749602965204744
$arr[7,37] = 1.749602965204744
$arr[7,38] = 1.749602965204744
$arr[7,39] = 1.749602965204744
$arr[7,40] = 1.749602965204744
$arr[7,41] = 1.749602965204744
$arr[7,42] = 1.749602965204744
$arr[7,43] = 1.749602965204744
$arr[7,44] = 1.749602965204744
$arr[7,45] = 1.749602965204744
$arr[7,46] = 1.749602965204744
$arr[7,47] = 1.749602965204744
$arr[7,48] = 1.749602965204744
$arr[7,49] = 1.749602965204744
$arr[7,50] = 1.749602965204744
$arr[7,51] = 1.749602965204744
$arr[8,0] = 42369
$arr[8,11] = 0.8024032016000104
$arr[8,12] = 0.8909614188480353
$arr[8,13] = 0.922773818606859
$arr[8,14] = 1.310895847186577
$arr[8,15] = 1.346932828201242
$arr[8,16] = 1.513781691628258
$arr[8,17] = 1.513781691628258
$arr[8,18] = 1.513781691628258
$arr[8,19] = 1.513781691628258
$arr[8,20] = 1.513781691628258
$arr[8,21] = 1.513781691628258
$arr[8,22] = 1.513781691628258
$arr[8,23] = 1.513781691628258
$arr[8,24] = 1.513781691628258
$arr[8,25] = 1.513781691628258
$arr[8,26] = 1.513781691628258
$arr[8,27] = 1.513781691628258
$arr[8,28] = 1.513781691628258
$arr[8,29] = 1.513781691628258
$arr[8,30] = 1.513781691628258
$arr[8,31] = 1.513781691628258
$arr[8,32] = 1.513781691628258
$arr[8,33] = 1.513781691628258
$arr[8,34] = 1.513781691628258
$arr[8,35] = 1.513781691628258
$arr[8,36] = 1.513781691628258
$arr[8,37] = 1.513781691628258
$arr[8,38] = 1.513781691628258
$arr[8,39] = 1.513781691628258
$arr[8,40] = 1.513781691628258
$arr[8,41] = 1.513781691628258
$arr[8,42] = 1.513781691628258
$arr[8,43] = 1.513781691628258
$arr[8,44] = 1.513781691628258
$arr[8,45] = 1.513781691628258
$arr[8,46] = 1.513781691628258
$arr[8,47] = 1.513781691628258
$arr[8,48] = 1.513781691628258
$arr[8,49] = 1.513781691628258
$arr[8,50] = 1.513781691628258
$arr[8,51] = 1.513781691628258
$arr[9,0] = 42735
$arr[9,13] = 0.8507004532711138
$arr[9,14] = 1.17312580610518
$arr[9,15] = 1.364302026343633
$arr[9,16] = 1.862478303083726
$arr[9,17] = 1.745747589686109
$arr[9,18] = 1.72540577912379
$arr[9,19] = 1.72540577912379
$arr[9,20] = 1.72540577912379
$arr[9,21] = 1.72540577912379
$arr[9,22] = 1.72540577912379
$arr[9,23] = 1.72540577912379
$arr[9,24] = 1.72540577912379
$arr[9,25] = 1.72540577912379
$arr[9,26] = 1.72540577912379
$arr[9,27] = 1.72540577912379
$arr[9,28] = 1.72540577912379
$arr[9,29] = 1.72540577912379
$arr[9,30] = 1.72540577912379
$arr[9,31] = 1.72540577912379
$arr[9,32] = 1.72540577912379
$arr[9,33] = 1.72540577912379
$arr[9,34] = 1.72540577912379
$arr[9,35] = 1.72540577912379
$arr[9,36] = 1.72540577912379
$arr[9,37] = 1.72540577912379
$arr[9,38] = 1.72540577912379
$arr[9,39] = 1.72540577912379
$arr[9,40] = 1.72540577912379
$arr[9,41] = 1.72540577912379
$arr[9,42] = 1.72540577912379
$arr[9,43] = 1.72540577912379
$arr[9,44] = 1.72540577912379
$arr[9,45] = 1.72540577912379
$arr[9,46] = 1.72540577912379
$arr[9,47] = 1.72540577912379
$arr[9,48] = 1.72540577912379
$arr[9,49] = 1.72540577912379
$arr[9,50] = 1.72540577912379
$arr[9,51] = 1.72540577912379
$arr[10,0] = 43100
$arr[10,15] = 1.342708276326299
$arr[10,16] = 1.804078246438934
$arr[10,17] = 1.644798626926303
$arr[10,18] = 1.639776099317536
$arr[10,19] = 1.843649045891893
$arr[10,20] = 2.026192376700298
$arr[10,21] = 2.026192376700298
$arr[10,22] = 2.026192376700298
$arr[10,23] = 2.026192376700298
$arr[10,24] = 2.026192376700298
$arr[10,25] = 2.026192376700298
$arr[10,26] = 2.026192376700298
$arr[10,27] = 2.026192376700298
$arr[10,28] = 2.026192376700298
$arr[10,29] = 2.026192376700298
$arr[10,30] = 2.026192376700298
$arr[10,31] = 2.026192376700298
$arr[10,32] = 2.026192376700298
$arr[10,33] = 2.026192376700298
$arr[10,34] = 2.026192376700298
$arr[10,35] = 2.026192376700298
$arr[10,36] = 2.026192376700298
$arr[10,37] = 2.026192376700298
$arr[10,38] = 2.026192376700298
$arr[10,39] = 2.026192376700298
$arr[10,40] = 2.026192376700298
$arr[10,41] = 2.026192376700298
$arr[10,42] = 2.026192376700298
$arr[10,43] = 2.026192376700298
$arr[10,44] = 2.026192376700298
$arr[10,45] = 2.026192376700298
$arr[10,46] = 2.026192376700298
$arr[10,47] = 2.026192376700298
$arr[10,48] = 2.026192376700298
$arr[10,49] = 2.026192376700298
$arr[10,50] = 2.026192376700298
$arr[10,51] = 2.026192376700298
$arr[11,0] = 43465
$arr[11,17] = 1.609625625599986
$arr[11,18] = 1.601529483008668
$arr[11,19] = 1.741128155516525
$arr[11,20] = 2.181728312936415
$arr[11,21] = 2.284406789710336
$arr[11,22] = 2.463589365374652
$arr[11,23] = 2.349806433215029
$arr[11,24] = 2.344166347125687
$arr[11,25] = 2.344166347125687
$arr[11,26] = 2.344166347125687
$arr[11,27] = 2.344166347125687
$arr[11,28] = 2.344166347125687
$arr[11,29] = 2.344166347125687
$arr[11,30] = 2.344166347125687
$arr[11,31] = 2.344166347125687
$arr[11,32] = 2.344166347125687
$arr[11,33] = 2.344166347125687
$arr[11,34] = 2.344166347125687
$arr[11,35] = 2.344166347125687
$arr[11,36] = 2.344166347125687
$arr[11,37] = 2.344166347125687
$arr[11,38] = 2.344166347125687
$arr[11,39] = 2.344166347125687
$arr[11,40] = 2.344166347125687
$arr[11,41] = 2.344166347125687
$arr[11,42] = 2.344166347125687
$arr[11,43] = 2.344166347125687
$arr[11,44] = 2.344166347125687
$arr[11,45] = 2.344166347125687
$arr[11,46] = 2.344166347125687
$arr[11,47] = 2.344166347125687
$arr[11,48] = 2.344166347125687
$arr[11,49] = 2.344166347125687
$arr[11,50] = 2.344166347125687
$arr[11,51] = 2.344166347125687
$arr[12,0] = 43830
$arr[12,19] = 1.723022434657207
$arr[12,20] = 1.99880983009828
$arr[12,21] = 1.990690441067144
$arr[12,22] = 2.149194501693219
$arr[12,23] = 2.036910005299108
$arr[12,24] = 2.010025322622599
$arr[12,25] = 1.665971362160357
$arr[12,26] = 1.332860091726285
$arr[12,27] = 1.029194292875912
$arr[12,28] = 0.9005461608770915
$arr[12,29] = 0.9005461608770915
$arr[12,30] = 0.9005461608770915
$arr[12,31] = 0.9005461608770915
$arr[12,32] = 0.9005461608770915
$arr[12,33] = 0.9005461608770915
$arr[12,34] = 0.9005461608770915
$arr[12,35] = 0.9005461608770915
$arr[12,36] = 0.9005461608770915
$arr[12,37] = 0.9005461608770915
$arr[12,38] = 0.9005461608770915
$arr[12,39] = 0.9005461608770915
$arr[12,40] = 0.9005461608770915
$arr[12,41] = 0.9005461608770915
$arr[12,42] = 0.9005461608770915
$arr[12,43] = 0.9005461608770915
$arr[12,44] = 0.9005461608770915
$arr[12,45] = 0.9005461608770915
$arr[12,46] = 0.9005461608770915
$arr[12,47] = 0.9005461608770915
$arr[12,48] = 0.9005461608770915
$arr[12,49] = 0.9005461608770915
$arr[12,50] = 0.9005461608770915
$arr[12,51] = 0.9005461608770915
$arr[13,0] = 44196
$arr[13,22] = 2.104406008906734
$arr[13,23] = 2.104406008906734
$arr[13,24] = 2.088155108730527
$arr[13,25] = 2.031292234149706
$arr[13,26] = 1.799885362733189
$arr[13,27] = 1.31420459445093
$arr[13,28] = 0.7771393814490102
$arr[13,29] = 0.4126128934655471
$arr[13,30] = 0.2336391425753925
$arr[13,31] = -4.43626840667447
$arr[13,32] = -4.43626840667447
$arr[13,33] = -4.43626840667447
$arr[13,34] = -4.43626840667447
$arr[13,35] = -4.43626840667447
$arr[13,36] = -4.43626840667447
$arr[13,37] = -4.43626840667447
$arr[13,38] = -4.43626840667447
$arr[13,39] = -4.43626840667447
$arr[13,40] = -4.43626840667447
$arr[13,41] = -4.43626840667447
$arr[13,42] = -4.43626840667447
$arr[13,43] = -4.43626840667447
$arr[13,44] = -4.43626840667447
$arr[13,45] = -4.43626840667447
$arr[13,46] = -4.43626840667447
$arr[13,47] = -4.43626840667447
$arr[13,48] = -4.43626840667447
$arr[13,49] = -4.43626840667447
$arr[13,50] = -4.43626840667447
$arr[13,51] = -4.43626840667447
$arr[14,0] = 44561
$arr[14,26] = 1.893295577996756
$arr[14,27] = 1.650113596657588
$arr[14,28] = 1.407367165006201
$arr[14,29] = 1.156986202028509
$arr[14,30] = 0.9207450904090253
$arr[14,31] = -2.63419394755392
$arr[14,32] = -2.71887004062904
$arr[14,33] = -2.96879819115512
$arr[14,34] = -2.010709456685855
$arr[14,35] = -1.513408827666285
$arr[14,36] = -1.513408827666285
$arr[14,37] = -1.513408827666285
$arr[14,38] = -1.513408827666285
$arr[14,39] = -1.513408827666285
$arr[14,40] = -1.513408827666285
$arr[14,41] = -1.513408827666285
$arr[14,42] = -1.513408827666285
$arr[14,43] = -1.513408827666285
$arr[14,44] = -1.513408827666285
$arr[14,45] = -1.513408827666285
$arr[14,46] = -1.513408827666285
$arr[14,47] = -1.513408827666285
$arr[14,48] = -1.513408827666285
$arr[14,49] = -1.513408827666285
$arr[14,50] = -1.513408827666285
$arr[14,51] = -1.513408827666285
$arr[15,0] = 44926
$arr[15,30] = 0.9554040357173665
$arr[15,31] = -1.6615457433243
$arr[15,32] = -1.937757788996253
$arr[15,33] = -2.438555173006141
$arr[15,34] = -1.14257141002756
$arr[15,35] = 0.7106578563214505
$arr[15,36] = 0.4582698374457683
$arr[15,37] = 1.154413086110817
$arr[15,38] = 1.5286818008164
$arr[15,39] = 1.618732201786743
$arr[15,40] = 1.618732201786743
$arr[15,41] = 1.618732201786743
$arr[15,42] = 1.618732201786743
$arr[15,43] = 1.618732201786743
$arr[15,44] = 1.618732201786743
$arr[15,45] = 1.618732201786743
$arr[15,46] = 1.618732201786743
$arr[15,47] = 1.618732201786743
$arr[15,48] = 1.618732201786743
$arr[15,49] = 1.618732201786743
$arr[15,50] = 1.618732201786743
$arr[15,51] = 1.618732201786743
$arr[16,0] = 45291
$arr[16,33] = -2.413633125962611
$arr[16,34] = -1.976964192800379
$arr[16,35] = -1.416208084524317
$arr[16,36] = -1.590220060268321
$arr[16,37] = -1.166698219025086
$arr[16,38] = -0.7118141543333012
$arr[16,39] = -0.4630595634534385
$arr[16,40] = 1.314675624401973
$arr[16,41] = 0.006126408955742235
$arr[16,42] = 0.001079933351455509
$arr[16,43] = -0.09609276733164585
$arr[16,44] = -0.09609276733164585
$arr[16,45] = -0.09609276733164585
$arr[16,46] = -0.09609276733164585
$arr[16,47] = -0.09609276733164585
$arr[16,48] = -0.09609276733164585
$arr[16,49] = -0.09609276733164585
$arr[16,50] = -0.09609276733164585
$arr[16,51] = -0.09609276733164585
$arr[17,0] = 45657
$arr[17,37] = -1.366617752737886
$arr[17,38] = -1.247796181364325
$arr[17,39] = -1.271569079498103
$arr[17,40] = 3.40905661301254
$arr[17,41] = 0.5370151562237302
$arr[17,42] = 0.6889047703476203
$arr[17,43] = 0.5285660612534882
$arr[17,44] = 0.1199358335146838
$arr[17,45] = -0.1634698065940632
$arr[17,46] = -0.00209793826797533
$arr[17,47] = -0.02761034355766023
$arr[17,48] = -0.02761034355766023
$arr[17,49] = -0.02761034355766023
$arr[17,50] = -0.02761034355766023
$arr[17,51] = -0.02761034355766023
$arr[18,0] = 46022
$arr[18,41] = 0.6286476167952193
$arr[18,42] = 0.8380511040219529
$arr[18,43] = 1.051341880573431
$arr[18,44] = 0.541919990776929
$arr[18,45] = -0.1145111565623136
$arr[18,46] = 0.286657616500996
$arr[18,47] = 0.08174908622293753
$arr[18,48] = 0.2676745853112728
$arr[18,49] = 0.4725905789402463
$arr[18,50] = 0.501314651583451
$arr[18,51] = 0.5152269879013183
$arr[19,0] = 46387
$arr[19,45] = -0.1001442164906607
$arr[19,46] = 0.2192992773568214
$arr[19,47] = 0.01413548204556303
$arr[19,48] = -0.05216314813395462
$arr[19,49] = 0.08564335928031852
$arr[19,50] = 0.2155158706220295
$arr[19,51] = 0.3332251551730891
$arr[20,0] = 46752
$arr[20,49] = 0.02760436504196662
$arr[20,50] = 0.1625720590646029
$arr[20,51] = 0.2843066506847514
$arr[21,0] = 47118

$ws.Range("A1:AZ22").Value2 = $arr

$wb.Save()
